$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.166.44'
$ws.Range("E2").Value = '  +3.67%  '
$ws.Range("D3").Value = '1.724.87'
$ws.Range("E3").Value = '  +2.64%  '
$ws.Range("E4").Value = '  -0.15%  '
$c = $ws.Range("D5")
$c.Style = "Normal"
$c.NumberFormat = "@"
$c.Value = '219.27'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.84%  '
$c = $ws.Range("D6")
$c.Style = "Normal"
$c.NumberFormat = "@"
$c.Value = '0.522'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.68%  '
$c = $ws.Range("D7")
$c.Style = "Normal"
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.16%  '
$c = $ws.Range("D8")
$c.Style = "Normal"
$c.NumberFormat = "@"
$c.Value = '24.35'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +13.93%  '
$ws.Range("E9").Value = '  +3.32%  '
$c = $ws.Range("D10")
$c.Style = "Normal"
$c.NumberFormat = "@"
$c.Value = '0.0634'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +1.71%  '
$ws.Range("E11").Value = '  +1.47%  '
$ws.Range("D12").Value = '1.967.52'
$ws.Range("E12").Value = '  +2.60%  '
$ws.Range("D13").Value = '1.728.49'
$ws.Range("E13").Value = '  +3.10%  '
$c = $ws.Range("D14")
$c.Style = "Normal"
$c.NumberFormat = "@"
$c.Value = '4.28'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +3.28%  '
$ws.Range("E15").Value = '  +4.68%  '
$c = $ws.Range("D16")
$c.Style = "Normal"
$c.NumberFormat = "@"
$c.Value = '67.57'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +1.98%  '
$ws.Range("D17").Value = '28.111.97'
$ws.Range("E17").Value = '  +3.49%  '
$c = $ws.Range("D18")
$c.Style = "Normal"
$c.NumberFormat = "@"
$c.Value = '243.76'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +1.79%  '
$c = $ws.Range("D19")
$c.Style = "Normal"
$c.NumberFormat = "@"
$c.Value = '8.02'
$c.Style = "Normal"
$ws.Range("D20").Value = '0.0₃0755'
$ws.Range("E20").Value = '  +1.60%  '
$ws.Range("E21").Value = '  -0.11%  '
$c = $ws.Range("D22")
$c.Style = "Normal"
$c.NumberFormat = "@"
$c.Value = '4.63'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +2.35%  '
$c = $ws.Range("D23")
$c.Style = "Normal"
$c.NumberFormat = "@"
$c.Value = '9.67'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +2.24%  '
$ws.Range("E24").Value = '  -0.23%  '
$c = $ws.Range("D25")
$c.Style = "Normal"
$c.NumberFormat = "@"
$c.Value = '149.28'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +1.49%  '
$c = $ws.Range("D26")
$c.Style = "Normal"
$c.NumberFormat = "@"
$c.Value = '7.51'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +3.66%  '
$c = $ws.Range("D27")
$c.Style = "Normal"
$c.NumberFormat = "@"
$c.Value = '16.71'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +2.26%  '
$ws.Range("E28").Value = '  +0.90%  '
$ws.Range("E29").Value = '  -0.22%  '
$ws.Range("E30").Value = '  +2.21%  '
$ws.Range("E31").Value = '  +1.59%  '
$ws.Range("E32").Value = '  +2.27%  '
$ws.Range("D33").Value = '1.498.67'
$ws.Range("E33").Value = '  -4.21%  '
$ws.Range("E34").Value = '  +1.43%  '
$ws.Range("E35").Value = '  -1.32%  '
$c = $ws.Range("D36")
$c.Style = "Normal"
$c.NumberFormat = "@"
$c.Value = '0.962'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +3.12%  '
$c = $ws.Range("D37")
$c.Style = "Normal"
$c.NumberFormat = "@"
$c.Value = '0.609'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +1.22%  '
$ws.Range("E38").Value = '  +0.54%  '
$c = $ws.Range("D39")
$c.Style = "Normal"
$c.NumberFormat = "@"
$c.Value = '0.0176'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +0.70%  '
$ws.Range("E40").Value = '  +1.51%  '
$c = $ws.Range("D41")
$c.Style = "Normal"
$c.NumberFormat = "@"
$c.Value = '70.85'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +2.26%  '
$ws.Range("E42").Value = '  +3.10%  '
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("E44").Value = '  +1.62%  '
$ws.Range("D45").Value = '1.872.63'
$ws.Range("E45").Value = '  +2.40%  '
$c = $ws.Range("D46")
$c.Style = "Normal"
$c.NumberFormat = "@"
$c.Value = '0.806'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +3.19%  '
$c = $ws.Range("D47")
$c.Style = "Normal"
$c.NumberFormat = "@"
$c.Value = '1.77'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +10.82%  '
$c = $ws.Range("D48")
$c.Style = "Normal"
$c.NumberFormat = "@"
$c.Value = '90.83'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +0.14%  '
$ws.Range("E49").Value = '  +5.92%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range("D50")
$c.Style = "Normal"
$c.NumberFormat = "@"
$c.Value = '0.105'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.74%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D51")
$c.Style = "Normal"
$c.NumberFormat = "@"
$c.Value = '8.21'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +0.84%  '
